$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary")

# --- Row 3 (Alan Mohamad) updates ---
$ws.Range("F3").Value = "Area (Aggressive)"
$ws.Range("G3").Value = 14717
$ws.Range("J3").Value = 114.44

# --- Row 4 (Ariel Young) updates ---
$ws.Range("C4").Value = "A17210559"
$ws.Range("F4").Value = "Area (Aggressive)"
$ws.Range("G4").Value = 14717
$ws.Range("J4").Value = 114.44

# --- New note row ---
$ws.Range("F6").Value = "Note: we also changed additional compiler settings (fitter effort, etc). Have included QPF files for reference."

# --- View settings ---
$excel.ActiveWindow.Zoom = 160
$ws.Range("F9").Select() | Out-Null
